$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Finalizacion del ejercicio 13: rename the four evaluation rows to reflect
# the finished/associated project names, and highlight the "13 Repaso..." row
# (C5:D5) with a new accent fill to mark it as finished. ---

$ws.Range("C5").Value = $ws.Range("C5").Value2 + " (EvaluacionBUG)"
$ws.Range("C6").Value = $ws.Range("C6").Value2 + " (JSEBEvaluación)"
$ws.Range("C7").Value = $ws.Range("C7").Value2 + " (Evaluacion1BryanUmanaGomez)"
$ws.Range("C8").Value = $ws.Range("C8").Value2 + " JavaBryanUmanaGomez)"

# Highlight the finished row (C5:D5) with a new accent color (theme Accent5).
$ws.Range("C5:D5").Interior.ThemeColor = 9

# Widen column C so the longer text fits.
$ws.Columns("C").ColumnWidth = 77.65

# Move the active selection, matching the saved view state.
$ws.Range("E6").Select()
